# Update cryptos list values (price and 1h volume change) per data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "71.676.39"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +2.76%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "4.032.31"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "525.63"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "148.58"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.628"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.19%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.741"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +1.83%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0000344"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "46.71"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +8.86%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "10.89"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +3.19%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "4.673.84"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +2.05%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "4.056.18"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +2.90%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "21.45"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +7.84%  "
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E19").Value = "  -1.51%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "71.657.16"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +2.75%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "444.30"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "3.61"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +5.78%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "94.91"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +7.15%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "14.44"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "12.17"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "4.04"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "11.12"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "37.26"
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "13.75"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +2.91%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "699.90"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("E32").Value = "  +1.85%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "6.98"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +15.12%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "67.98"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.0₃0912"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").Value = "  +1.19%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "41.49"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("E38").Value = "  +3.45%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "3.53"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +17.27%  "
$ws.Range("E40").Value = "  +0.12%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0495"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  -0.11%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.87"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  +1.25%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.54"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("E47").Value = "  -2.36%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "9.32"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +7.19%  "
$ws.Range("E49").Value = "  +17.34%  "
$ws.Range("E50").Value = "  +1.10%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0₆0345"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -5.86%  "

Write-Output "Updated cryptos list values"
